# Auto-generated edit script: updates crypto price/volume table
# to match the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force text storage so numeric-looking strings (e.g. "57.50")
    # are not coerced into Number cells and lose formatting,
    # then drop back to the Normal style so no stray number format
    # is left applied to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "36.451.85"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3
$ws.Range("D3").Value = "1.943.12"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue $ws.Range("D5") "243.28"
$ws.Range("E5").Value = "  -0.24%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.615"
$ws.Range("E6").Value = "  -0.94%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "57.50"
$ws.Range("E8").Value = "  -1.21%  "

# Row 9
$ws.Range("E9").Value = "  -2.46%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0844"
$ws.Range("E10").Value = "  -1.47%  "

# Row 11
$ws.Range("E11").Value = "  -0.90%  "

# Row 12
$ws.Range("D12").Value = "2.229.35"
$ws.Range("E12").Value = "  -1.04%  "

# Row 13
Set-TextValue $ws.Range("D13") "21.33"
$ws.Range("E13").Value = "  -4.35%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.812"
$ws.Range("E14").Value = "  -2.78%  "

# Row 15
Set-TextValue $ws.Range("D15") "13.46"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16
Set-TextValue $ws.Range("D16") "5.17"
$ws.Range("E16").Value = "  -3.50%  "

# Row 17
$ws.Range("D17").Value = "1.955.89"
$ws.Range("E17").Value = "  -0.68%  "

# Row 18
$ws.Range("D18").Value = "36.431.18"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
Set-TextValue $ws.Range("D19") "69.33"
$ws.Range("E19").Value = "  -2.73%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("E20").Value = "  -2.92%  "

# Row 21
Set-TextValue $ws.Range("D21") "228.48"
$ws.Range("E21").Value = "  -1.23%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.99"
$ws.Range("E22").Value = "  -3.03%  "

# Row 23
$ws.Range("E23").Value = "  -0.18%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.35"
$ws.Range("E24").Value = "  -7.01%  "

# Row 25
$ws.Range("E25").Value = "  +1.29%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.19"
$ws.Range("E26").Value = "  -4.23%  "

# Row 27
Set-TextValue $ws.Range("D27") "161.70"
$ws.Range("E27").Value = "  -3.17%  "

# Row 28
$ws.Range("E28").Value = "  +9.25%  "

# Row 29
Set-TextValue $ws.Range("D29") "19.17"
$ws.Range("E29").Value = "  -4.16%  "

# Row 30
$ws.Range("E30").Value = "  -1.20%  "

# Row 31
$ws.Range("E31").Value = "  -4.69%  "

# Row 32
$ws.Range("E32").Value = "  -3.45%  "

# Row 33
$ws.Range("E33").Value = "  -3.81%  "

# Row 34
$ws.Range("E34").Value = "  -4.21%  "

# Row 35
Set-TextValue $ws.Range("D35") "6.22"
$ws.Range("E35").Value = "  +4.31%  "

# Row 36
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("E37").Value = "  -1.15%  "

# Row 38
$ws.Range("E38").Value = "  +0.84%  "

# Row 39
$ws.Range("E39").Value = "  +8.19%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0980"
$ws.Range("E40").Value = "  +1.97%  "

# Row 41
$ws.Range("E41").Value = "  +0.40%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0209"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D43") "1.15"
$ws.Range("E43").Value = "  -3.16%  "

# Row 44
Set-TextValue $ws.Range("D44") "15.99"
$ws.Range("E44").Value = "  +1.48%  "

# Row 45
$ws.Range("D45").Value = "1.343.42"
$ws.Range("E45").Value = "  -0.45%  "

# Row 46
$ws.Range("E46").Value = "  -2.96%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "86.57"
$ws.Range("E47").Value = "  -2.61%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D48") "7.20"
$ws.Range("E48").Value = "  -1.22%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.82"
$ws.Range("E49").Value = "  +0.05%  "

# Row 50
$ws.Range("D50").Value = "2.120.00"
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
Set-TextValue $ws.Range("D51") "43.31"
$ws.Range("E51").Value = "  -3.81%  "
